$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from row 2 through row 301 (all rows with value 45186 -> 45188)
$range = $ws.Range("C2:C301")
$range.Value = 45188
